$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply header styling (bold, centered, bordered) to the new AB:AN block ---
# Row 1 header band (style matches existing A1/O1 band style)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AB1:AN1").PasteSpecial(-4122) | Out-Null

# Row 2 sub-header band
$ws.Range("A2").Copy() | Out-Null
$ws.Range("AB2:AN2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Write cell values (headers, labels and data) ---
$ws.Range("B1").Value = "msg_count_twitter"
$ws.Range("O1").Value = "msg_count_twitter_engage"
$ws.Range("AB1").Value = "msg_count_facebook"
$ws.Range("B2").Value = "sum"
$ws.Range("C2").Value = "mean"
$ws.Range("D2").Value = "std"
$ws.Range("E2").Value = "min"
$ws.Range("F2").Value = "q25"
$ws.Range("G2").Value = "median"
$ws.Range("H2").Value = "q75"
$ws.Range("I2").Value = "max"
$ws.Range("J2").Value = "count"
$ws.Range("K2").Value = "msg_per_mus"
$ws.Range("L2").Value = "active_mus_n"
$ws.Range("M2").Value = "active_mus_pc"
$ws.Range("N2").Value = "active_mus_pc_z"
$ws.Range("O2").Value = "sum"
$ws.Range("P2").Value = "mean"
$ws.Range("Q2").Value = "std"
$ws.Range("R2").Value = "min"
$ws.Range("S2").Value = "q25"
$ws.Range("T2").Value = "median"
$ws.Range("U2").Value = "q75"
$ws.Range("V2").Value = "max"
$ws.Range("W2").Value = "count"
$ws.Range("X2").Value = "msg_per_mus"
$ws.Range("Y2").Value = "active_mus_n"
$ws.Range("Z2").Value = "active_mus_pc"
$ws.Range("AA2").Value = "active_mus_pc_z"
$ws.Range("AB2").Value = "sum"
$ws.Range("AC2").Value = "mean"
$ws.Range("AD2").Value = "std"
$ws.Range("AE2").Value = "min"
$ws.Range("AF2").Value = "q25"
$ws.Range("AG2").Value = "median"
$ws.Range("AH2").Value = "q75"
$ws.Range("AI2").Value = "max"
$ws.Range("AJ2").Value = "count"
$ws.Range("AK2").Value = "msg_per_mus"
$ws.Range("AL2").Value = "active_mus_n"
$ws.Range("AM2").Value = "active_mus_pc"
$ws.Range("AN2").Value = "active_mus_pc_z"
$ws.Range("A3").Value = "size"
$ws.Range("A4").Value = "huge"
$ws.Range("B4").Value = 32564
$ws.Range("C4").Value = 2713.7
$ws.Range("D4").Value = 3185.9
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1835.5
$ws.Range("H4").Value = 4087.2
$ws.Range("I4").Value = 10793
$ws.Range("J4").Value = 12
$ws.Range("K4").Value = 4070.5
$ws.Range("L4").Value = 8
$ws.Range("M4").Value = 66.7
$ws.Range("N4").Value = 0.2
$ws.Range("O4").Value = 57730
$ws.Range("P4").Value = 4810.8
$ws.Range("Q4").Value = 5758.1
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 52.5
$ws.Range("T4").Value = 2818
$ws.Range("U4").Value = 7519.2
$ws.Range("V4").Value = 18016
$ws.Range("W4").Value = 12
$ws.Range("X4").Value = 6414.4
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 75
$ws.Range("AA4").Value = 0.6
$ws.Range("AB4").Value = 11708
$ws.Range("AC4").Value = 975.7
$ws.Range("AD4").Value = 918.8
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 218.2
$ws.Range("AG4").Value = 992
$ws.Range("AH4").Value = 1503.5
$ws.Range("AI4").Value = 2991
$ws.Range("AJ4").Value = 12
$ws.Range("AK4").Value = 1300.9
$ws.Range("AL4").Value = 9
$ws.Range("AM4").Value = 75
$ws.Range("AN4").Value = 1.3
$ws.Range("A5").Value = "large"
$ws.Range("B5").Value = 760187
$ws.Range("C5").Value = 1526.5
$ws.Range("D5").Value = 2393.3
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 142
$ws.Range("G5").Value = 825
$ws.Range("H5").Value = 1918.2
$ws.Range("I5").Value = 19900
$ws.Range("J5").Value = 498
$ws.Range("K5").Value = 1881.7
$ws.Range("L5").Value = 404
$ws.Range("M5").Value = 81.09999999999999
$ws.Range("N5").Value = 1.4
$ws.Range("O5").Value = 339665
$ws.Range("P5").Value = 682.1
$ws.Range("Q5").Value = 1780.5
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 35.8
$ws.Range("T5").Value = 205
$ws.Range("U5").Value = 597.5
$ws.Range("V5").Value = 20111
$ws.Range("W5").Value = 498
$ws.Range("X5").Value = 824.4
$ws.Range("Y5").Value = 412
$ws.Range("Z5").Value = 82.7
$ws.Range("AA5").Value = 1.2
$ws.Range("AB5").Value = 267908
$ws.Range("AC5").Value = 538
$ws.Range("AD5").Value = 592
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0
$ws.Range("AG5").Value = 425.5
$ws.Range("AH5").Value = 844.8
$ws.Range("AI5").Value = 4241
$ws.Range("AJ5").Value = 498
$ws.Range("AK5").Value = 832
$ws.Range("AL5").Value = 322
$ws.Range("AM5").Value = 64.7
$ws.Range("AN5").Value = 0.8
$ws.Range("A6").Value = "medium"
$ws.Range("B6").Value = 674674
$ws.Range("C6").Value = 780.9
$ws.Range("D6").Value = 2020.4
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 191
$ws.Range("H6").Value = 965.8
$ws.Range("I6").Value = 43003
$ws.Range("J6").Value = 864
$ws.Range("K6").Value = 1159.2
$ws.Range("L6").Value = 582
$ws.Range("M6").Value = 67.40000000000001
$ws.Range("N6").Value = 0.3
$ws.Range("O6").Value = 213128
$ws.Range("P6").Value = 246.7
$ws.Range("Q6").Value = 1703.7
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 28
$ws.Range("U6").Value = 162.5
$ws.Range("V6").Value = 46992
$ws.Range("W6").Value = 864
$ws.Range("X6").Value = 358.8
$ws.Range("Y6").Value = 594
$ws.Range("Z6").Value = 68.8
$ws.Range("AA6").Value = 0.2
$ws.Range("AB6").Value = 385317
$ws.Range("AC6").Value = 446
$ws.Range("AD6").Value = 1552.3
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 0
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 548.5
$ws.Range("AI6").Value = 36948
$ws.Range("AJ6").Value = 864
$ws.Range("AK6").Value = 1044.2
$ws.Range("AL6").Value = 369
$ws.Range("AM6").Value = 42.7
$ws.Range("AN6").Value = -0.2
$ws.Range("A7").Value = "small"
$ws.Range("B7").Value = 813150
$ws.Range("C7").Value = 445.6
$ws.Range("D7").Value = 1507.5
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 5
$ws.Range("H7").Value = 363
$ws.Range("I7").Value = 47580
$ws.Range("J7").Value = 1825
$ws.Range("K7").Value = 865.1
$ws.Range("L7").Value = 940
$ws.Range("M7").Value = 51.5
$ws.Range("N7").Value = -1
$ws.Range("O7").Value = 255885
$ws.Range("P7").Value = 140.2
$ws.Range("Q7").Value = 615.1
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 1
$ws.Range("U7").Value = 56
$ws.Range("V7").Value = 11227
$ws.Range("W7").Value = 1825
$ws.Range("X7").Value = 268.8
$ws.Range("Y7").Value = 952
$ws.Range("Z7").Value = 52.2
$ws.Range("AA7").Value = -0.9
$ws.Range("AB7").Value = 445848
$ws.Range("AC7").Value = 244.3
$ws.Range("AD7").Value = 1197.5
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 0
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 23
$ws.Range("AI7").Value = 36858
$ws.Range("AJ7").Value = 1825
$ws.Range("AK7").Value = 958.8
$ws.Range("AL7").Value = 465
$ws.Range("AM7").Value = 25.5
$ws.Range("AN7").Value = -1
$ws.Range("A8").Value = "unknown"
$ws.Range("B8").Value = 123254
$ws.Range("C8").Value = 880.4
$ws.Range("D8").Value = 4123.7
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 2.5
$ws.Range("H8").Value = 772.5
$ws.Range("I8").Value = 47571
$ws.Range("J8").Value = 140
$ws.Range("K8").Value = 1688.4
$ws.Range("L8").Value = 73
$ws.Range("M8").Value = 52.1
$ws.Range("N8").Value = -0.9
$ws.Range("O8").Value = 42654
$ws.Range("P8").Value = 304.7
$ws.Range("Q8").Value = 1808.6
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 0.5
$ws.Range("U8").Value = 105
$ws.Range("V8").Value = 20884
$ws.Range("W8").Value = 140
$ws.Range("X8").Value = 609.3
$ws.Range("Y8").Value = 70
$ws.Range("Z8").Value = 50
$ws.Range("AA8").Value = -1.1
$ws.Range("AB8").Value = 75835
$ws.Range("AC8").Value = 541.7
$ws.Range("AD8").Value = 2079.8
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 0
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 239.8
$ws.Range("AI8").Value = 20246
$ws.Range("AJ8").Value = 140
$ws.Range("AK8").Value = 1849.6
$ws.Range("AL8").Value = 41
$ws.Range("AM8").Value = 29.3
$ws.Range("AN8").Value = -0.8

# --- Merge the new header band cell, mirroring B1:N1 / O1:AA1 ---
$ws.Range("AB1:AN1").Merge() | Out-Null
